# edit.ps1 - Applies the "Export PDFs locally (project root)" edits to
# Completion Report.docx via Word COM interop (iron_native runtime).

$d = $word.ActiveDocument

# Helper: Find & replace plain text (affects ALL matches in the story because
# this runtime's Find.Execute with Wrap=wdFindContinue(1) replaces every
# occurrence it finds while scanning to the end of the story).
function Replace-Text($findText, $replaceText) {
    $null = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $replaceText, 2)
}

# Helper: replace the run(s) covering an exact piece of visible text with a
# single freshly built <w:r> whose inner XML is supplied by the caller. Used
# for structural run-level tweaks (adding/removing <w:lastRenderedPageBreak/>)
# that plain text Find/Replace cannot express.
function Set-RunXml($searchText, $innerXml) {
    $probe = $d.Content
    $found = $probe.Find.Execute($searchText)
    if (-not $found) {
        throw "Set-RunXml: text not found: $searchText"
    }
    $target = $d.Range($probe.Start, $probe.End)
    $payload = '<?xml version="1.0" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p><w:r>' + $innerXml + '</w:r></w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($payload)
}

# ---------------------------------------------------------------------------
# 1) "Project Status:" line and the six "Version 1 Final Status" bullets all
#    drop their leading "✅ " emoji run (the emoji + following space run are
#    removed, leaving just the plain text behind).
# ---------------------------------------------------------------------------
Replace-Text "✅ " ""

# ---------------------------------------------------------------------------
# 2) "Issue: Camtel (ISP) connectivity..." - drop the spell-check wrapper
#    around "Camtel" by collapsing the surrounding runs into one.
# ---------------------------------------------------------------------------
Replace-Text ": Camtel (ISP) connectivity problems prevented timely reporting" `
             ": Camtel (ISP) connectivity problems prevented timely reporting"

# ---------------------------------------------------------------------------
# 3) "Resolution: Continued development offline, finalized documentation
#    upon connectivity restoration" paragraph gets substantially rewritten.
#    Locate it precisely (it's the "Resolution" run right after "Camtel").
# ---------------------------------------------------------------------------
$afterCamtel = $d.Content
$null = $afterCamtel.Find.Execute("Camtel")
$tail = $d.Range($afterCamtel.End, $d.Content.End)
$null = $tail.Find.Execute("Resolution")
$tail.Text = "Improvision (not resolution)"
$tail.Collapse(0)

$restStart = $tail.Start
$rest = $d.Range($restStart, $d.Content.End)
$null = $rest.Find.Execute(": Continued development offline, finalized documentation upon connectivity restoration")
$newTail = ": Continued development offline, finalized documentation upon personal " + `
           "mobile data. But was still slowed down really bad. The offline state makes it difficult " + `
           "me write mostly what I can remember from the optimization " + [char]0x2013 + " The server " + `
           "has been inaccessible remote. It has been that way for more than a week."
$rest.Text = $newTail

# ---------------------------------------------------------------------------
# 4) Move <w:lastRenderedPageBreak/> from the "Integration with external
#    laboratory systems" run to the "Enhanced patient portal capabilities"
#    run (the paragraph immediately before it).
# ---------------------------------------------------------------------------
Set-RunXml "Enhanced patient portal capabilities" `
    "<w:lastRenderedPageBreak/><w:t>Enhanced patient portal capabilities</w:t>"
Set-RunXml "Integration with external laboratory systems" `
    "<w:t>Integration with external laboratory systems</w:t>"

# ---------------------------------------------------------------------------
# 5) Sign-off block: replace placeholder names with the real author.
# ---------------------------------------------------------------------------
Replace-Text " Cyber-Lord Development Team" " Nzenong K. Mc Braxton Development Team"
Replace-Text " [Your Name]" " Nzenong K. Mc Braxton"
